$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (fecha 45009) is inserted as row 220; all rows
# that were previously 220..305 shift down to 221..306 (dimension grows to
# A1:R306). Insert a whole row so everything below shifts automatically.
$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new record's data.
$ws.Cells.Item(220, 1).Value = 8
$ws.Cells.Item(220, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(220, 3).Value = "Coquimbo"
$ws.Cells.Item(220, 4).Value = 45009
$ws.Cells.Item(220, 5).Value = 4
$ws.Cells.Item(220, 6).Value = 100112037
$ws.Cells.Item(220, 7).Value = "Cebollín"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 1300
$ws.Cells.Item(220, 11).Value = 1000
$ws.Cells.Item(220, 12).Value = 1200
$ws.Cells.Item(220, 13).Value = 1100
$ws.Cells.Item(220, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(220, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(220, 16).Value = 183
$ws.Cells.Item(220, 17).Value = 6
$ws.Cells.Item(220, 18).Value = "Hortaliza"
